$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers in I1 and J1, matching style of existing headers (H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$lastRow = 28
for ($r = 2; $r -le $lastRow; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
